$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Countries")

# Fix the misspelled "Addres" header -> "address"
$ws.Range("C1").Value = "address"

# Reflect the resulting selection change observed after the edit
$ws.Activate()
$ws.Range("E11").Select()
